# [PHOENIX-5860] updated Collect fee trade License feature
#
# - update the legacy trade license number on the "tradeLocationDetails" sheet
# - swap the flammables sub-category example on the "tradeDetails" sheet
# - add a new "legacyDetails" sheet (Collect-fee / legacy-trade lookup data)
# - leave each sheet's selection where the author last left it, with the new
#   sheet becoming the active tab

$wb = $excel.ActiveWorkbook

# --- tradeOwnerDetails: no data change, just leave the cursor on A2 ---
$ws1 = $wb.Worksheets.Item("tradeOwnerDetails")
[void]$ws1.Range("A2").Select()

# --- tradeLocationDetails: 1016017647 -> 1016047857 ---
$ws2 = $wb.Worksheets.Item("tradeLocationDetails")
$ws2.Range("B2").Value = "1016047857"
[void]$ws2.Range("C2").Select()

# --- tradeDetails: "Petrol or Diesel bunks" -> "Acetylene Gas" ---
$ws3 = $wb.Worksheets.Item("tradeDetails")
$ws3.Range("E2").Value = "Acetylene Gas"
[void]$ws3.Range("B2").Select()

# --- new sheet: legacyDetails (appended after tradeDetails) ---
$newSheet = $wb.Worksheets.Add($null, $ws3)
$newSheet.Name = "legacyDetails"
$newSheet.Range("A1").Value = "data Name"
$newSheet.Range("B1").Value = "Fee Details"
$newSheet.Range("A2").Value = "legency Trade"
[void]$newSheet.Range("B11").Select()
